# Fruta / hortaliza, semanal
# Insert 7 new weekly price rows for Durazno (Macroferia Regional de Talca)
# above the existing last rows (which simply shift down from 508-511 to 515-518).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 7 new rows just above the current last data rows (508..511).
$ws.Rows("508:514").Insert()

$rows = @(
    @{ Row = 508; K = "Carson";       L = "Especial";                M = 210; N = 17000; O = 17000; P = 17000; Q = "$/bandeja 15 kilos granel";   S = 1133 },
    @{ Row = 509; K = "Carson";       L = "Extra (doble especial)";  M = 150; N = 19000; O = 19000; P = 19000; Q = "$/bandeja 15 kilos granel";   S = 1267 },
    @{ Row = 510; K = "Carson";       L = "Primera";                 M = 280; N = 15000; O = 15000; P = 15000; Q = "$/bandeja 15 kilos granel";   S = 1000 },
    @{ Row = 511; K = "Royal Glory";  L = "Especial";                M = 180; N = 17000; O = 17000; P = 17000; Q = "$/bandeja 15 kilos granel";   S = 1133 },
    @{ Row = 512; K = "Royal Glory";  L = "Primera";                 M = 210; N = 15000; O = 15000; P = 15000; Q = "$/bandeja 15 kilos granel";   S = 1000 },
    @{ Row = 513; K = "Toscana";      L = "Especial";                M = 220; N = 17000; O = 17000; P = 17000; Q = "$/bandeja 15 kilos granel";   S = 1133 },
    @{ Row = 514; K = "Toscana";      L = "Extra (doble especial)";  M = 240; N = 19000; O = 19000; P = 19000; Q = "$/bandeja 15 kilos granel";   S = 1267 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 5
    $ws.Range("B$row").Value = "Macroferia Regional de Talca"
    $ws.Range("C$row").Value = "Maule"
    $ws.Range("D$row").Value = 44939
    $ws.Range("E$row").Value = 7
    $ws.Range("F$row").Value = "Fruta"
    $ws.Range("G$row").Value = 100103
    $ws.Range("H$row").Value = "Frutos de hueso (carozo)"
    $ws.Range("I$row").Value = 100103004
    $ws.Range("J$row").Value = "Durazno"
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = "Región de O'Higgins"
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = 15
}
